{"js": "// The \"\u0421\u0432\u043e\u0439\u0441\u0442\u0432\u0430 \u043f\u0441\u0435\u0432\u0434\u043e\u043a\u043e\u043c\u043f\u043e\u043d\u0435\u043d\u0442\u043e\u0432\" heading paragraph and its three bullet\n// items (numId=2 list) are removed and collapsed into a single trailing\n// empty \"List Paragraph\"-styled paragraph that just carries a \"_GoBack\"\n// bookmark (this is what Word leaves behind after the last edit position\n// once the content that used to live there has been deleted).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph index 6 is \"\u0421\u0432\u043e\u0439\u0441\u0442\u0432\u0430 \u043f\u0441\u0435\u0432\u0434\u043e\u043a\u043e\u043c\u043f\u043e\u043d\u0435\u043d\u0442\u043e\u0432\" and 7-9 are the three\n// bullet points that follow it - delete all four (reverse order so the\n// indices of paragraphs still to be removed stay valid).\nconst firstIndex = 6;\nconst lastIndex = 9;\nfor (let i = lastIndex; i >= firstIndex; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// Add the replacement empty paragraph at the end of the body, then use\n// insertOoxml to give it the exact final shape (pStyle \"a3\" with Arial\n// paragraph-mark formatting, no list numbering, plus the _GoBack bookmark)\n// without leaving a stray empty run behind.\nconst newParagraph = body.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:pStyle w:val=\"a3\"/>' +\n  '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/></w:rPr>' +\n  '</w:pPr>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nnewParagraph.getRange().insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The \"\u0421\u0432\u043e\u0439\u0441\u0442\u0432\u0430 \u043f\u0441\u0435\u0432\u0434\u043e\u043a\u043e\u043c\u043f\u043e\u043d\u0435\u043d\u0442\u043e\u0432\" heading paragraph and its three bullet\n# items (the numId=2 list) are removed and collapsed into a single trailing\n# empty \"List Paragraph\"-styled paragraph that just carries a \"_GoBack\"\n# bookmark (this is what Word leaves behind after the last edit position\n# once the content that used to live there has been deleted).\n\n$d = $word.ActiveDocument\n\n# Paragraph 7 (1-indexed COM) is \"\u0421\u0432\u043e\u0439\u0441\u0442\u0432\u0430 \u043f\u0441\u0435\u0432\u0434\u043e\u043a\u043e\u043c\u043f\u043e\u043d\u0435\u043d\u0442\u043e\u0432\"; delete\n# everything from its start through the very end of the document (this\n# removes paragraphs 7-10, i.e. the heading plus its three bullet points).\n$firstPara = $d.Paragraphs(7)\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$deleteRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)\n$deleteRange.Delete()\n\n# Add the replacement paragraph at the end of the body. A second, throwaway\n# paragraph is added right after it so the first one can be fully replaced\n# (paragraph mark included) via InsertXML without losing the final\n# paragraph mark of the document.\n$newPara = $d.Paragraphs.Add()\n$helperPara = $d.Paragraphs.Add()\n\n$targetIndex = $d.Paragraphs.Count - 1\n$targetPara = $d.Paragraphs($targetIndex)\n$targetRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End)\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:pStyle w:val=\"a3\"/>' +\n  '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/></w:rPr>' +\n  '</w:pPr>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$targetRange.InsertXML($ooxml)\n\n# Remove the now-trailing throwaway paragraph, leaving our replacement\n# paragraph as the last paragraph in the document.\n$finalCount = $d.Paragraphs.Count\n$trailingRange = $d.Range($d.Paragraphs($finalCount - 1).Range.End, $d.Paragraphs($finalCount).Range.End)\n$trailingRange.Delete()\n"}
